$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "1.00", "7.80", "10.50").
# Force text format on those cells first so Excel keeps the exact string
# instead of silently re-parsing them as numbers and dropping trailing zeros.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.163.24"
$ws.Range("E2").Value = "  +4.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.330.59"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.75"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.93"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.354.81"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +5.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.71"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.752.91"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.936.94"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.363.84"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.70"
$ws.Range("E21").Value = "  +5.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.06"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("E25").Value = "  +8.11%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.951"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.80"
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("E28").Value = "  +9.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.56"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +6.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.29"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.932"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  +5.60%  "
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.85"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +5.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.34"
$ws.Range("E43").Value = "  +4.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.26"
$ws.Range("E44").Value = "  +10.31%  "
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0937"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0505"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.565"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0219"
$ws.Range("E49").Value = "  +5.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.379"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  +11.18%  "
